$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 64
$ws.Range("H64").Value = 4010.3125
$ws.Range("I64").Value = 3733.2
$ws.Range("J64").Value = 5000
$ws.Range("K64").Value = 3733.2
$ws.Range("L64").Value = 5000
$ws.Range("M64").Value = -3485.2
$ws.Range("N64").Value = -5496
# row 67
$ws.Range("H67").Value = 4010.3125
$ws.Range("I67").Value = 3733.2
$ws.Range("J67").Value = 5000
$ws.Range("K67").Value = 3733.2
$ws.Range("L67").Value = 5000
$ws.Range("M67").Value = -2875.2
$ws.Range("N67").Value = -6716
# row 87
$ws.Range("H87").Value = 13288.535
$ws.Range("J87").Value = 13288.535
$ws.Range("L87").Value = 13288.535
$ws.Range("N87").Value = -15784.535
# row 90
$ws.Range("H90").Value = 13288.535
$ws.Range("J90").Value = 13288.535
$ws.Range("L90").Value = 39865.605
$ws.Range("N90").Value = -52345.605
# row 113
$ws.Range("H113").Value = 5631.317
$ws.Range("I113").Value = 2257.1538
$ws.Range("J113").Value = 11479.866
$ws.Range("K113").Value = 2257.1538
$ws.Range("L113").Value = 11479.866
$ws.Range("M113").Value = 996.8462
$ws.Range("N113").Value = -17987.866
# row 132
$ws.Range("H132").Value = 68115.7
$ws.Range("I132").Value = 96436.71000000001
$ws.Range("J132").Value = 2033.3334
$ws.Range("K132").Value = 289310.13
$ws.Range("L132").Value = 6100.0002
$ws.Range("M132").Value = -286780.13
$ws.Range("N132").Value = -11160.0002

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 2
$ws.Range("H2").Value = 1531.3914
$ws.Range("I2").Value = 1434.5555
$ws.Range("J2").Value = 1880
$ws.Range("K2").Value = 1434.5555
$ws.Range("L2").Value = 1880
$ws.Range("M2").Value = -1321.5555
$ws.Range("N2").Value = -2106
# row 7
$ws.Range("H7").Value = 27200
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 27200
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 27200
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -27428
# row 116
$ws.Range("H116").Value = 1531.3914
$ws.Range("I116").Value = 1434.5555
$ws.Range("J116").Value = 1880
$ws.Range("K116").Value = 1434.5555
$ws.Range("L116").Value = 1880
$ws.Range("M116").Value = 859.4445000000001
$ws.Range("N116").Value = -6468
# row 122
$ws.Range("H122").Value = 2433.04
$ws.Range("I122").Value = 2190.111
$ws.Range("J122").Value = 3057.7144
$ws.Range("K122").Value = 6570.333
$ws.Range("L122").Value = 9173.143199999999
$ws.Range("M122").Value = -4120.333
$ws.Range("N122").Value = -14073.1432
# row 132
$ws.Range("H132").Value = 2288.919
$ws.Range("I132").Value = 2039.8572
$ws.Range("K132").Value = 6119.571599999999
$ws.Range("M132").Value = -3589.571599999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 3
$ws.Range("H3").Value = 1531.3914
$ws.Range("I3").Value = 1434.5555
$ws.Range("J3").Value = 1880
$ws.Range("K3").Value = 1434.5555
$ws.Range("L3").Value = 1880
$ws.Range("M3").Value = -1320.5555
$ws.Range("N3").Value = -2108
# row 94
$ws.Range("H94").Value = 581.63635
$ws.Range("I94").Value = 796.3333
$ws.Range("J94").Value = 501.125
$ws.Range("K94").Value = 796.3333
$ws.Range("L94").Value = 501.125
$ws.Range("M94").Value = -345.3333
$ws.Range("N94").Value = -1403.125

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 22
$ws.Range("H22").Value = 5292639.5
$ws.Range("I22").Value = 9261906
$ws.Range("K22").Value = 9261906
$ws.Range("M22").Value = -9261556
# row 58
$ws.Range("H58").Value = 2220.3416
$ws.Range("I58").Value = 1501.05
$ws.Range("J58").Value = 2905.3809
$ws.Range("K58").Value = 1501.05
$ws.Range("L58").Value = 2905.3809
$ws.Range("M58").Value = -1298.05
$ws.Range("N58").Value = -3311.3809
# row 99
$ws.Range("H99").Value = 18843.666
$ws.Range("I99").Value = 2015.5
$ws.Range("K99").Value = 2015.5
$ws.Range("M99").Value = -517.5
# row 126
$ws.Range("H126").Value = 18843.666
$ws.Range("I126").Value = 2015.5
$ws.Range("K126").Value = 6046.5
$ws.Range("M126").Value = -3576.5
# row 136
$ws.Range("H136").Value = 2220.3416
$ws.Range("I136").Value = 1501.05
$ws.Range("J136").Value = 2905.3809
$ws.Range("K136").Value = 4503.15
$ws.Range("L136").Value = 8716.1427
$ws.Range("M136").Value = -1953.15
$ws.Range("N136").Value = -13816.1427

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 107
$ws.Range("H107").Value = 504685.06
$ws.Range("J107").Value = 851300.7
$ws.Range("L107").Value = 2553902.1
$ws.Range("N107").Value = -2557742.1
# row 131
$ws.Range("H131").Value = 846.03
$ws.Range("J131").Value = 868.7659
$ws.Range("L131").Value = 2606.2977
$ws.Range("N131").Value = -12686.2977

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 70
$ws.Range("H70").Value = 5753.2354
$ws.Range("I70").Value = 4884.4
$ws.Range("J70").Value = 6439.1577
$ws.Range("K70").Value = 4884.4
$ws.Range("L70").Value = 6439.1577
$ws.Range("M70").Value = -4614.4
$ws.Range("N70").Value = -6979.1577
# row 73
$ws.Range("H73").Value = 5753.2354
$ws.Range("I73").Value = 4884.4
$ws.Range("J73").Value = 6439.1577
$ws.Range("K73").Value = 4884.4
$ws.Range("L73").Value = 6439.1577
$ws.Range("M73").Value = -3948.4
$ws.Range("N73").Value = -8311.1577
# row 97
$ws.Range("H97").Value = 2358.7144
$ws.Range("I97").Value = 2583.3333
$ws.Range("J97").Value = 1011
$ws.Range("K97").Value = 2583.3333
$ws.Range("L97").Value = 1011
$ws.Range("M97").Value = -2087.3333
$ws.Range("N97").Value = -2003
# row 102
$ws.Range("H102").Value = 5498019
$ws.Range("I102").Value = 9618782
$ws.Range("J102").Value = 3669
$ws.Range("K102").Value = 9618782
$ws.Range("L102").Value = 3669
$ws.Range("M102").Value = -9617160
$ws.Range("N102").Value = -6913
# row 122
$ws.Range("H122").Value = 2913.742
$ws.Range("I122").Value = 2345.9048
$ws.Range("K122").Value = 7037.714399999999
$ws.Range("M122").Value = -4587.714399999999
# row 126
$ws.Range("H126").Value = 2052.8865
$ws.Range("I126").Value = 1844.8928
$ws.Range("J126").Value = 2416.875
$ws.Range("K126").Value = 5534.678400000001
$ws.Range("L126").Value = 7250.625
$ws.Range("M126").Value = -3064.678400000001
$ws.Range("N126").Value = -12190.625
# row 132
$ws.Range("H132").Value = 2301.639
$ws.Range("I132").Value = 1966.15
$ws.Range("J132").Value = 2721
$ws.Range("K132").Value = 5898.450000000001
$ws.Range("L132").Value = 8163
$ws.Range("M132").Value = -3368.450000000001
$ws.Range("N132").Value = -13223
# row 136
$ws.Range("H136").Value = 9645.553
$ws.Range("J136").Value = 9645.553
$ws.Range("L136").Value = 28936.659
$ws.Range("N136").Value = -34036.659

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 7
$ws.Range("H7").Value = 44570.918
$ws.Range("I7").Value = 57952.945
$ws.Range("J7").Value = 4424.8335
$ws.Range("K7").Value = 57952.945
$ws.Range("L7").Value = 4424.8335
$ws.Range("M7").Value = -57840.945
$ws.Range("N7").Value = -4648.8335
# row 40
$ws.Range("H40").Value = 23292.834
$ws.Range("I40").Value = 32957.688
$ws.Range("K40").Value = 32957.688
$ws.Range("M40").Value = -32821.688
# row 122
$ws.Range("H122").Value = 8549230
$ws.Range("I122").Value = 18519802
$ws.Range("J122").Value = 3025
$ws.Range("K122").Value = 55559406
$ws.Range("L122").Value = 9075
$ws.Range("M122").Value = -55556956
$ws.Range("N122").Value = -13975
# row 126
$ws.Range("H126").Value = 44570.918
$ws.Range("I126").Value = 57952.945
$ws.Range("J126").Value = 4424.8335
$ws.Range("K126").Value = 173858.835
$ws.Range("L126").Value = 13274.5005
$ws.Range("M126").Value = -171388.835
$ws.Range("N126").Value = -18214.5005
# row 132
$ws.Range("H132").Value = 8990.326999999999
$ws.Range("I132").Value = 8201.236999999999
$ws.Range("J132").Value = 10489.6
$ws.Range("K132").Value = 24603.711
$ws.Range("L132").Value = 31468.8
$ws.Range("M132").Value = -22073.711
$ws.Range("N132").Value = -36528.8

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 46
$ws.Range("H46").Value = 41085.734
$ws.Range("J46").Value = 41085.734
$ws.Range("L46").Value = 41085.734
$ws.Range("N46").Value = -41547.734
# row 122
$ws.Range("H122").Value = 73982.07000000001
$ws.Range("I122").Value = 93085.82000000001
$ws.Range("K122").Value = 279257.46
$ws.Range("M122").Value = -276807.46
# row 126
$ws.Range("H126").Value = 28716.084
$ws.Range("I126").Value = 38945.348
$ws.Range("J126").Value = 2120
$ws.Range("K126").Value = 116836.044
$ws.Range("L126").Value = 6360
$ws.Range("M126").Value = -114366.044
$ws.Range("N126").Value = -11300
# row 127
$ws.Range("H127").Value = 36750
$ws.Range("J127").Value = 36750
$ws.Range("L127").Value = 36750
$ws.Range("N127").Value = -46670
# row 132
$ws.Range("H132").Value = 1211.6538
$ws.Range("I132").Value = 842.4103
$ws.Range("J132").Value = 2319.3845
$ws.Range("K132").Value = 2527.2309
$ws.Range("L132").Value = 6958.1535
$ws.Range("M132").Value = 2.76909999999998
$ws.Range("N132").Value = -12018.1535
# row 134
$ws.Range("H134").Value = 41085.734
$ws.Range("J134").Value = 41085.734
$ws.Range("L134").Value = 123257.202
$ws.Range("N134").Value = -128327.202
